{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list (under the \"Impact\"\n// sub-heading) into four short, impact-focused accomplishment statements,\n// dropping the two testimony / FEC-compliance \"job duty\" bullets entirely.\n//\n// NOTE: several of the old bullet strings here are near-duplicates of text\n// that legitimately lives elsewhere in the resume (e.g. under \"Partner -\n// Siege Analytics\"), so this script locates the target paragraphs\n// structurally (by walking forward from the \"KEY ACHIEVEMENTS AND IMPACT\"\n// heading) rather than doing a blind whole-document text search/replace.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" section heading.\nlet sectionIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionIdx = i;\n    break;\n  }\n}\nif (sectionIdx === -1) {\n  throw new Error('Could not find the \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\n\n// Collect the bullet paragraphs (\"\u2022 ...\") that follow it, stopping at the\n// next non-bullet paragraph (the following section heading).\nconst bulletParas = [];\nfor (let i = sectionIdx + 1; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t.startsWith(\"\\u2022\")) {\n    bulletParas.push(items[i]);\n  } else if (bulletParas.length > 0) {\n    break;\n  }\n}\n\nif (bulletParas.length !== 6) {\n  throw new Error(\n    `Expected 6 existing bullets under KEY ACHIEVEMENTS AND IMPACT, found ${bulletParas.length}`\n  );\n}\n\nconst newText = [\n  \"\\u2022 Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard\",\n  \"\\u2022 Reduced polling margins from \\u00b14.2% to \\u00b12.1%\",\n  \"\\u2022 Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\",\n  \"\\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\",\n];\n\n// Rewrite the first four bullets in place.\nfor (let i = 0; i < newText.length; i++) {\n  bulletParas[i].insertText(newText[i], \"Replace\");\n}\n\n// Remove the trailing two bullets (testimony / FEC compliance) entirely.\nbulletParas[4].delete();\nbulletParas[5].delete();\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list (under the \"Impact\"\n# sub-heading) into four short, impact-focused accomplishment statements,\n# dropping the two testimony / FEC-compliance \"job duty\" bullets entirely.\n#\n# NOTE: several of the old bullet strings here are near-duplicates of text\n# that legitimately lives elsewhere in the resume (e.g. under \"Partner -\n# Siege Analytics\"), so this script locates the target paragraphs\n# structurally (by walking forward from the \"KEY ACHIEVEMENTS AND IMPACT\"\n# heading) rather than doing a blind whole-document text search/replace.\n\n$d = $word.ActiveDocument\n\n# Find the \"KEY ACHIEVEMENTS AND IMPACT\" heading paragraph.\n$sectionIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionIdx = $i\n        break\n    }\n}\nif ($sectionIdx -eq -1) {\n    throw \"Could not find the 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# Collect the bullet paragraphs (\"<bullet> ...\") that follow it, stopping at\n# the next non-bullet paragraph (the following section heading).\n$bullet = [char]0x2022\n$bulletIdx = New-Object System.Collections.ArrayList\nfor ($i = $sectionIdx + 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t.StartsWith($bullet)) {\n        [void]$bulletIdx.Add($i)\n    } elseif ($bulletIdx.Count -gt 0) {\n        break\n    }\n}\n\nif ($bulletIdx.Count -ne 6) {\n    throw \"Expected 6 existing bullets under KEY ACHIEVEMENTS AND IMPACT, found \" + $bulletIdx.Count\n}\n\n$pm = [char]0x00B1\n\n$newTexts = @(\n    ($bullet + \" Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard\"),\n    ($bullet + \" Reduced polling margins from \" + $pm + \"4.2% to \" + $pm + \"2.1%\"),\n    ($bullet + \" Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\"),\n    ($bullet + \" Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\")\n)\n\n# Rewrite the first four bullets in place.\nfor ($k = 0; $k -lt 4; $k++) {\n    $idx = [int]$bulletIdx[$k]\n    $d.Paragraphs.Item($idx).Range.Text = $newTexts[$k]\n}\n\n# Remove the trailing two bullets (testimony / FEC compliance) entirely.\n# Delete the higher index first so the lower index stays valid.\n$d.Paragraphs.Item([int]$bulletIdx[5]).Range.Delete()\n$d.Paragraphs.Item([int]$bulletIdx[4]).Range.Delete()\n"}
